# Insert a new data row at row 165 (shifting existing rows 165..230 down to 166..231)
# then populate the new row with the weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(165).Insert()

$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 45134
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = 100112043
$ws.Cells.Item(165, 7).Value = "Pepino ensalada"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 100
$ws.Cells.Item(165, 11).Value = 10000
$ws.Cells.Item(165, 12).Value = 11000
$ws.Cells.Item(165, 13).Value = 10500
$ws.Cells.Item(165, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(165, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(165, 16).Value = 175
$ws.Cells.Item(165, 17).Value = 60
$ws.Cells.Item(165, 18).Value = "Hortaliza"
